# Append 3 new trade rows (rows 8-10) to "WorkSheet 1", matching the
# BIIBNamedTrade20.xlsx diff: new data rows with Principle, Start Principle,
# BuyPrice, SellPrice, IsShortSell, Price Change %, Date, Profitable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Cells.Item(8, 1).Value = 9972.19
$ws.Cells.Item(8, 2).Value = 9906.81
$ws.Cells.Item(8, 3).Value = 305.24
$ws.Cells.Item(8, 4).Value = 307.24
$ws.Cells.Item(8, 5).Value = $false
$ws.Cells.Item(8, 6).Value = 0.66
$ws.Cells.Item(8, 7).Value = 42613.766770833332
$ws.Cells.Item(8, 8).Value = $true

# Row 9
$ws.Cells.Item(9, 1).Value = 9990.14
$ws.Cells.Item(9, 2).Value = 9972.19
$ws.Cells.Item(9, 3).Value = 307.68
$ws.Cells.Item(9, 4).Value = 308.24
$ws.Cells.Item(9, 5).Value = $false
$ws.Cells.Item(9, 6).Value = 0.18
$ws.Cells.Item(9, 7).Value = 42614.674803240741
$ws.Cells.Item(9, 8).Value = $true

# Row 10
$ws.Cells.Item(10, 1).Value = 9986.14
$ws.Cells.Item(10, 2).Value = 9990.14
$ws.Cells.Item(10, 3).Value = 307.95999999999998
$ws.Cells.Item(10, 4).Value = 307.83
$ws.Cells.Item(10, 5).Value = $false
$ws.Cells.Item(10, 6).Value = -0.04
$ws.Cells.Item(10, 7).Value = 42615.751863425925
$ws.Cells.Item(10, 8).Value = $false

# Apply the same date style (style index 1 / numFmtId 22) used by the rest of
# column G by copying the format from the cell directly above the new block.
$ws.Range("G7").Copy()
$ws.Range("G8:G10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
